# Restauración desde backup y correcciones de restablecer contraseña
#
# Summary of changes applied:
#  1. Insert a new sheet "COMENTARIOS_GRAFICAS" right after "GRAFICAS" and
#     populate it with a small parameter table (TIPO_COMENTARIO / Free /
#     Basico / Pro / Premium).
#  2. Append a new "GRAFICAS_INTELIGENTES" row to the "PLANES" sheet and
#     bold the label column for that sheet.
#  3. Tidy up row-height overrides on "EN_ANALISIS".
#  4. Update the remembered selections on a couple of sheets and move the
#     active tab to "PLANES".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet: COMENTARIOS_GRAFICAS (inserted after GRAFICAS)
# ---------------------------------------------------------------------
$graficas = $wb.Worksheets.Item("GRAFICAS")
$comentarios = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $graficas)
$comentarios.Name = "COMENTARIOS_GRAFICAS"
$comentarios.Columns("A").ColumnWidth = 20.14

# Write the new label strings in this exact order so the shared-string
# table ends up with the same ordering as the target workbook.
$comentarios.Range("A3").Value = "Comentario GPT"
$comentarios.Range("A1").Value = "TIPO_COMENTARIO"
$comentarios.Range("A2").Value = "Comentario Python"

$comentarios.Range("B1").Value = "Free"
$comentarios.Range("C1").Value = "Basico"
$comentarios.Range("D1").Value = "Pro"
$comentarios.Range("E1").Value = "Premium"

$comentarios.Range("B2").Value = 0
$comentarios.Range("C2").Value = 1
$comentarios.Range("D2").Value = 1
$comentarios.Range("E2").Value = 1

$comentarios.Range("B3").Value = 0
$comentarios.Range("C3").Value = 0
$comentarios.Range("D3").Value = 0
$comentarios.Range("E3").Value = 1

# Re-use the existing header styles from TARJETAS (style used for the
# first / remaining header cells) so no redundant style entries appear.
$tarjetas = $wb.Worksheets.Item("TARJETAS")
$tarjetas.Range("A1").Copy() | Out-Null
$comentarios.Range("A1").PasteSpecial(-4122) | Out-Null
$tarjetas.Range("E1").Copy() | Out-Null
$comentarios.Range("B1:E1").PasteSpecial(-4122) | Out-Null

$comentarios.Range("E3").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. PLANES: add GRAFICAS_INTELIGENTES row + bold label column
# ---------------------------------------------------------------------
$planes = $wb.Worksheets.Item("PLANES")

$planes.Range("A10").Value = "GRAFICAS_INTELIGENTES"
$planes.Range("B10").Value = 0
$planes.Range("C10").Value = 0
$planes.Range("D10").Value = 0
$planes.Range("E10").Value = 1

# Bold the parameter-name column (A2:A10), matching the new cellXf.
$planes.Range("A2:A10").Font.Bold = $true

$planes.Range("A11").Select() | Out-Null
$planes.Activate() | Out-Null

# ---------------------------------------------------------------------
# 3. EN_ANALISIS: clear the stale row-height overrides on rows 6 & 7
# ---------------------------------------------------------------------
$enAnalisis = $wb.Worksheets.Item("EN_ANALISIS")
$enAnalisis.Rows(6).AutoFit() | Out-Null
$enAnalisis.Rows(7).AutoFit() | Out-Null
$enAnalisis.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Remembered selections on the other sheets
# ---------------------------------------------------------------------
$kpis = $wb.Worksheets.Item("KPIS_FINANCIEROS")
$kpis.Range("A6:XFD6").Select() | Out-Null

$graficas.Range("A1:XFD3").Select() | Out-Null

# ---------------------------------------------------------------------
# Leave PLANES as the active sheet (matches workbookView activeTab).
# ---------------------------------------------------------------------
$planes.Activate() | Out-Null
